$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update metadata type for "sector-de-actividad" (column A) from dimension to measure
$ws.Range("A2").Value = "iaest-measure:sector-de-actividad"
$ws.Range("A3").Value = "medida"
$ws.Range("A4").Value = "xsd:int"

# Update metadata type for "sexo" (column F) from dimension to measure
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("F3").Value = "medida"
$ws.Range("F4").Value = "xsd:int"

# Remove the now-obsolete mapping-file row (row 5)
$ws.Range("A5:H5").Delete()
